$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "stok" (stock) column (D) to 10 for every existing product row
# (rows 2-37 already had a D value, rows 38-66 did not have one at all).
$ws.Range("D2:D66").Value = 10

# Add the two new products at the bottom of the table (rows 67 and 68),
# mirroring the pattern used by the existing rows (formulas in A/C derive
# the display name / generated code from columns E/F).
$ws.Range("E67").Value = "Sinocare BA801"
$ws.Range("F67").Value = "Tensimeter digital"
$ws.Range("B67").Value = 250000
$ws.Range("D67").Value = 10
$ws.Range("A67").Formula = '=F67 & " - " & E67'
$ws.Range("C67").Formula = '=UPPER(LEFT(F67,4) & "-" & LEFT(E67,2) & "-" & TEXT(ROW(A67)-1,"000"))'

$ws.Range("E68").Value = "Sinocare 50"
$ws.Range("F68").Value = "Alat Cek Gula Darah"
$ws.Range("B68").Value = 150000
$ws.Range("D68").Value = 10
$ws.Range("A68").Formula = '=F68 & " - " & E68'
$ws.Range("C68").Formula = '=UPPER(LEFT(F68,4) & "-" & LEFT(E68,2) & "-" & TEXT(ROW(A68)-1,"000"))'
